$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.734.89'
$ws.Range("E2").Value = '  -2.66%  '
$ws.Range("D3").Value = '1.780.06'
$ws.Range("E3").Value = '  -2.16%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5118'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3790'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07781'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.56%  '
$ws.Range("E11").Value = '  -2.47%  '
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.199'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.09'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.43%  '
$ws.Range("D15").Value = '1.776.16'
$ws.Range("E15").Value = '  -2.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.166'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.25'
$ws.Range("D17").Style = "Normal"
$ws.Range("E18").Value = '  -5.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06555'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.05%  '
$ws.Range("E22").Value = '  -2.74%  '
$ws.Range("D23").Value = '27.784.44'
$ws.Range("E23").Value = '  -2.63%  '
$ws.Range("E24").Value = '  -3.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.237'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.81%  '
$ws.Range("D28").Value = '1.982.46'
$ws.Range("E28").Value = '  -2.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.350'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1069'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.030'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.632'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.473'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07056'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02308'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.731'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.38%  '
$ws.Range("E38").Value = '  -5.19%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.019'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.48%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.48'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6076'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.150'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.320'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.04'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5921'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.713'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.70%  '
$ws.Range("E48").Value = '  +1.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.195'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.892'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06814'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.33%  '
